$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.896.31"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "1.741.81"
$ws.Range("E3").Value = "  -1.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5231"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2757"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "

$ws.Range("D11").Value = "1.737.33"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07116"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6448"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.521"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").Value = "25.889.23"
$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006666"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "

$ws.Range("D22").Value = "1.959.39"
$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.270"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.778"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.172"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.802"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08340"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.731"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.568"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04522"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.612"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9779"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6228"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.697"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01589"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.920"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9994"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3870"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7374"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.010"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05334"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1126"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.239"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.653"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.31%  "
